# "Evaluate values before parsing tests and add 'enabled' parameter to Delete tag"
#
# - Header row test cells: C1 now carries the new "disabled={{...}}" flavour
#   of the <<delete>> tag, and E1 is switched from the old formula-driven
#   OnlyValues/Delete text to the "disabled={{disableEColumnDeletion}}" text.
# - F1 was just a styled placeholder with no content - drop it entirely.
# - Row 1 no longer needs an explicit row height (let Excel auto-fit it).
# - Moves the saved selection from J9 to K9.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C1").Value = "<<delete disabled={{disableCColumnDeletion}}>>"
$ws.Range("E1").Value = "<<delete disabled={{disableEColumnDeletion}}>>"

# F1 was just a styled placeholder with no content - drop it entirely.
$ws.Range("F1").Clear()

# Let row 1 auto-size again instead of keeping the stale explicit height.
$ws.Rows.Item(1).AutoFit()

# Move the persisted selection from J9 to K9.
$null = $ws.Range("K9").Select()
